$wb = $excel.ActiveWorkbook

# "Repayment schedule" is the 3rd sheet (sheet3.xml)
$ws = $wb.Worksheets.Item(3)

# Insert a new blank column before column N (14th column), shifting the
# existing "Late"/"heading"/"Outstanding" columns one to the right, to make
# room for the new "Variable Instalments" related column.
$srcColumnWidth = $ws.Columns.Item(13).ColumnWidth
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $srcColumnWidth

# Make "Repayment schedule" the active/selected sheet with R12 selected.
$ws.Activate()
[void]$ws.Range("R12").Select()
